$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the values for columns A..T (20 columns) of one
# data row. Rows correspond to worksheet rows 2..13 (row 1 is the header).
$rowsData = @(
    @("ECs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 22.113303, 66.339909, 0.9168319126933797, 0.9168319126933796, 3, 1, 5.105980333333334, 15.317941, 0.04494800244635684, 0.04494800244635684, 112.910090223041, 1016.190812007369, 0.04120976305464005, 0.04120976305464005),
    @("ECs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 22.113303, 66.339909, 0.9168319126933797, 0.9168319126933796, 3, 1, 92.82552866666667, 278.476586, 0.8171441754985934, 0.8171441754985934, 2052.679041541186, 18474.11137387068, 0.7491838573686301, 0.74918385736863),
    @("ECs", "Pdgfd", "Pdgfrb", "M1", 3, 1, 22.113303, 66.339909, 0.9168319126933797, 0.9168319126933796, 1, 0.3333333333333333, 0.04351533333333333, 0.130546, 0.0003830659699865731, 0.0003830659699865732, 0.962267751146, 8.660409760314, 0.0003512071059505346, 0.0003512071059505347),
    @("ECs", "Pdgfd", "Pdgfrb", "sCs", 3, 1, 22.113303, 66.339909, 0.9168319126933797, 0.9168319126933796, 3, 1, 15.62246733333333, 46.867402, 0.1375247560850632, 0.1375247560850632, 345.464353749602, 3109.179183746418, 0.126087085164159, 0.126087085164159),
    @("FAPs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 1.846222666666667, 5.538667999999999, 0.0765455915264161, 0.07654559152641609, 3, 1, 5.105980333333334, 15.317941, 0.04494800244635684, 0.04494800244635684, 9.426776626954222, 84.840989642588, 0.003440571435187182, 0.003440571435187182),
    @("FAPs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 1.846222666666667, 5.538667999999999, 0.0765455915264161, 0.07654559152641609, 3, 1, 92.82552866666667, 278.476586, 0.8171441754985934, 0.8171441754985934, 171.3765950697164, 1542.389355627448, 0.0625487842759054, 0.06254878427590539),
    @("FAPs", "Pdgfd", "Pdgfrb", "M1", 3, 1, 1.846222666666667, 5.538667999999999, 0.0765455915264161, 0.07654559152641609, 1, 0.3333333333333333, 0.04351533333333333, 0.130546, 0.0003830659699865731, 0.0003830659699865732, 0.08033899474755554, 0.7230509527279999, 0.0000293220112662626, 0.0000293220112662626),
    @("FAPs", "Pdgfd", "Pdgfrb", "sCs", 3, 1, 1.846222666666667, 5.538667999999999, 0.0765455915264161, 0.07654559152641609, 3, 1, 15.62246733333333, 46.867402, 0.1375247560850632, 0.1375247560850632, 28.84255330005955, 259.582979700536, 0.01052691380405726, 0.01052691380405725),
    @("sCs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 0.1597296666666667, 0.479189, 0.006622495780204159, 0.006622495780204158, 3, 1, 5.105980333333334, 15.317941, 0.04494800244635684, 0.04494800244635684, 0.8155765366498889, 7.340188829849001, 0.0002976679565296044, 0.0002976679565296044),
    @("sCs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 0.1597296666666667, 0.479189, 0.006622495780204159, 0.006622495780204158, 3, 1, 92.82552866666667, 278.476586, 0.8171441754985934, 0.8171441754985934, 14.82699075208378, 133.442916768754, 0.005411533854057841, 0.005411533854057841),
    @("sCs", "Pdgfd", "Pdgfrb", "M1", 3, 1, 0.1597296666666667, 0.479189, 0.006622495780204159, 0.006622495780204158, 1, 0.3333333333333333, 0.04351533333333333, 0.130546, 0.0003830659699865731, 0.0003830659699865732, 0.006950689688222221, 0.062556207194, 0.000002536852769775894, 0.000002536852769775894),
    @("sCs", "Pdgfd", "Pdgfrb", "sCs", 3, 1, 0.1597296666666667, 0.479189, 0.006622495780204159, 0.006622495780204158, 3, 1, 15.62246733333333, 46.867402, 0.1375247560850632, 0.1375247560850632, 2.495371499664222, 22.458343496978, 0.0009107571168469373, 0.0009107571168469372)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $values = $rowsData[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
